$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = -19.25730036044673
$ws.Range("C2").Value2 = 1.956293383358081
$ws.Range("D2").Value2 = -19.25730036044673
$ws.Range("E2").Value2 = -19.25730036044673
$ws.Range("F2").Value2 = -19.25730036044673
$ws.Range("G2").Value2 = -19.25730036044673
$ws.Range("H2").Value2 = -19.25730036044673
$ws.Range("I2").Value2 = -19.25730036044673
$ws.Range("J2").Value2 = -19.25730036044673
$ws.Range("K2").Value2 = -19.25730036044673
$ws.Range("B3").Value2 = -19.25730036044673
$ws.Range("C3").Value2 = -19.25730036044673
$ws.Range("D3").Value2 = -19.25730036044673
$ws.Range("E3").Value2 = -19.25730036044673
$ws.Range("F3").Value2 = -19.25730036044673
$ws.Range("G3").Value2 = -19.25730036044673
$ws.Range("H3").Value2 = -19.25730036044673
$ws.Range("I3").Value2 = 1.323985745062338
$ws.Range("J3").Value2 = -19.25730036044673
$ws.Range("K3").Value2 = -19.25730036044673
$ws.Range("B4").Value2 = -19.25730036044673
$ws.Range("C4").Value2 = 1.979363803394709
$ws.Range("D4").Value2 = 1.694141599937667
$ws.Range("E4").Value2 = -19.25730036044673
$ws.Range("F4").Value2 = 3.478722667332117
$ws.Range("G4").Value2 = -19.25730036044673
$ws.Range("H4").Value2 = 1.674932417021585
$ws.Range("I4").Value2 = -19.25730036044673
$ws.Range("J4").Value2 = 1.28648008814529
$ws.Range("K4").Value2 = -19.25730036044673
$ws.Range("B5").Value2 = -19.25730036044673
$ws.Range("C5").Value2 = 1.686408479583733
$ws.Range("D5").Value2 = -19.25730036044673
$ws.Range("E5").Value2 = -19.25730036044673
$ws.Range("F5").Value2 = -19.25730036044673
$ws.Range("G5").Value2 = 3.040632671610326
$ws.Range("H5").Value2 = -19.25730036044673
$ws.Range("I5").Value2 = -19.25730036044673
$ws.Range("J5").Value2 = -19.25730036044673
$ws.Range("K5").Value2 = -19.25730036044673
$ws.Range("B6").Value2 = -19.25730036044673
$ws.Range("C6").Value2 = -19.25730036044673
$ws.Range("D6").Value2 = -19.25730036044673
$ws.Range("E6").Value2 = -19.25730036044673
$ws.Range("F6").Value2 = -19.25730036044673
$ws.Range("G6").Value2 = -19.25730036044673
$ws.Range("H6").Value2 = -19.25730036044673
$ws.Range("I6").Value2 = -19.25730036044673
$ws.Range("J6").Value2 = -19.25730036044673
$ws.Range("K6").Value2 = -19.25730036044673
$ws.Range("B7").Value2 = 2.459877764486468
$ws.Range("C7").Value2 = -19.25730036044673
$ws.Range("D7").Value2 = -19.25730036044673
$ws.Range("E7").Value2 = -19.25730036044673
$ws.Range("F7").Value2 = -19.25730036044673
$ws.Range("G7").Value2 = -19.25730036044673
$ws.Range("H7").Value2 = -19.25730036044673
$ws.Range("I7").Value2 = -19.25730036044673
$ws.Range("J7").Value2 = -19.25730036044673
$ws.Range("K7").Value2 = -19.25730036044673
$ws.Range("B8").Value2 = -19.25730036044673
$ws.Range("C8").Value2 = -19.25730036044673
$ws.Range("D8").Value2 = -19.25730036044673
$ws.Range("E8").Value2 = 1.820950147324369
$ws.Range("F8").Value2 = -19.25730036044673
$ws.Range("G8").Value2 = -19.25730036044673
$ws.Range("H8").Value2 = -19.25730036044673
$ws.Range("I8").Value2 = -19.25730036044673
$ws.Range("J8").Value2 = -19.25730036044673
$ws.Range("K8").Value2 = -19.25730036044673
$ws.Range("B9").Value2 = 3.857808873009537
$ws.Range("C9").Value2 = -19.25730036044673
$ws.Range("D9").Value2 = -19.25730036044673
$ws.Range("E9").Value2 = -19.25730036044673
$ws.Range("F9").Value2 = -19.25730036044673
$ws.Range("G9").Value2 = -19.25730036044673
$ws.Range("H9").Value2 = -19.25730036044673
$ws.Range("I9").Value2 = -19.25730036044673
$ws.Range("J9").Value2 = -19.25730036044673
$ws.Range("K9").Value2 = -19.25730036044673
$ws.Range("B10").Value2 = -19.25730036044673
$ws.Range("C10").Value2 = -19.25730036044673
$ws.Range("D10").Value2 = -19.25730036044673
$ws.Range("E10").Value2 = -19.25730036044673
$ws.Range("F10").Value2 = -19.25730036044673
$ws.Range("G10").Value2 = -19.25730036044673
$ws.Range("H10").Value2 = -19.25730036044673
$ws.Range("I10").Value2 = 1.826724886467073
$ws.Range("J10").Value2 = -19.25730036044673
$ws.Range("K10").Value2 = -19.25730036044673
$ws.Range("B11").Value2 = -19.25730036044673
$ws.Range("C11").Value2 = -19.25730036044673
$ws.Range("D11").Value2 = -19.25730036044673
$ws.Range("E11").Value2 = 2.928439237183959
$ws.Range("F11").Value2 = -19.25730036044673
$ws.Range("G11").Value2 = 2.679266054463187
$ws.Range("H11").Value2 = -19.25730036044673
$ws.Range("I11").Value2 = -19.25730036044673
$ws.Range("J11").Value2 = -19.25730036044673
$ws.Range("K11").Value2 = -19.25730036044673
$ws.Range("B12").Value2 = -19.25730036044673
$ws.Range("C12").Value2 = -19.25730036044673
$ws.Range("D12").Value2 = -19.25730036044673
$ws.Range("E12").Value2 = -19.25730036044673
$ws.Range("F12").Value2 = -19.25730036044673
$ws.Range("G12").Value2 = -19.25730036044673
$ws.Range("H12").Value2 = -19.25730036044673
$ws.Range("I12").Value2 = -19.25730036044673
$ws.Range("J12").Value2 = -19.25730036044673
$ws.Range("K12").Value2 = -19.25730036044673
$ws.Range("B13").Value2 = -19.25730036044673
$ws.Range("C13").Value2 = -19.25730036044673
$ws.Range("D13").Value2 = -19.25730036044673
$ws.Range("E13").Value2 = 2.494770298792863
$ws.Range("F13").Value2 = -19.25730036044673
$ws.Range("G13").Value2 = -19.25730036044673
$ws.Range("H13").Value2 = -19.25730036044673
$ws.Range("I13").Value2 = -19.25730036044673
$ws.Range("J13").Value2 = 1.86387203829443
$ws.Range("K13").Value2 = -19.25730036044673
$ws.Range("B14").Value2 = -19.25730036044673
$ws.Range("C14").Value2 = -19.25730036044673
$ws.Range("D14").Value2 = 1.512938466566365
$ws.Range("E14").Value2 = -19.25730036044673
$ws.Range("F14").Value2 = -19.25730036044673
$ws.Range("G14").Value2 = -19.25730036044673
$ws.Range("H14").Value2 = -19.25730036044673
$ws.Range("I14").Value2 = -19.25730036044673
$ws.Range("J14").Value2 = -19.25730036044673
$ws.Range("K14").Value2 = -19.25730036044673
$ws.Range("B15").Value2 = -19.25730036044673
$ws.Range("C15").Value2 = -19.25730036044673
$ws.Range("D15").Value2 = 1.728700657504364
$ws.Range("E15").Value2 = -19.25730036044673
$ws.Range("F15").Value2 = -19.25730036044673
$ws.Range("G15").Value2 = -19.25730036044673
$ws.Range("H15").Value2 = -19.25730036044673
$ws.Range("I15").Value2 = -19.25730036044673
$ws.Range("J15").Value2 = -19.25730036044673
$ws.Range("K15").Value2 = -19.25730036044673
$ws.Range("B16").Value2 = -19.25730036044673
$ws.Range("C16").Value2 = -19.25730036044673
$ws.Range("D16").Value2 = -19.25730036044673
$ws.Range("E16").Value2 = -19.25730036044673
$ws.Range("F16").Value2 = -19.25730036044673
$ws.Range("G16").Value2 = -19.25730036044673
$ws.Range("H16").Value2 = -19.25730036044673
$ws.Range("I16").Value2 = -19.25730036044673
$ws.Range("J16").Value2 = 2.044232116695251
$ws.Range("K16").Value2 = -19.25730036044673
$ws.Range("B17").Value2 = -19.25730036044673
$ws.Range("C17").Value2 = 2.15570200923836
$ws.Range("D17").Value2 = 1.857088245981997
$ws.Range("E17").Value2 = -19.25730036044673
$ws.Range("F17").Value2 = -19.25730036044673
$ws.Range("G17").Value2 = -19.25730036044673
$ws.Range("H17").Value2 = 2.0013697544314
$ws.Range("I17").Value2 = 1.960653197120457
$ws.Range("J17").Value2 = 2.307485963628553
$ws.Range("K17").Value2 = -19.25730036044673
$ws.Range("B18").Value2 = -19.25730036044673
$ws.Range("C18").Value2 = -19.25730036044673
$ws.Range("D18").Value2 = -19.25730036044673
$ws.Range("E18").Value2 = -19.25730036044673
$ws.Range("F18").Value2 = -19.25730036044673
$ws.Range("G18").Value2 = -19.25730036044673
$ws.Range("H18").Value2 = 2.021741681997814
$ws.Range("I18").Value2 = 2.088003629485321
$ws.Range("J18").Value2 = 2.276820581419839
$ws.Range("K18").Value2 = -19.25730036044673
$ws.Range("B19").Value2 = -19.25730036044673
$ws.Range("C19").Value2 = -19.25730036044673
$ws.Range("D19").Value2 = 2.037641725790754
$ws.Range("E19").Value2 = -19.25730036044673
$ws.Range("F19").Value2 = -19.25730036044673
$ws.Range("G19").Value2 = -19.25730036044673
$ws.Range("H19").Value2 = 1.604044817439888
$ws.Range("I19").Value2 = 1.766327445078502
$ws.Range("J19").Value2 = -19.25730036044673
$ws.Range("K19").Value2 = -19.25730036044673
$ws.Range("B20").Value2 = -19.25730036044673
$ws.Range("C20").Value2 = 1.0142048421077
$ws.Range("D20").Value2 = 1.519763077160697
$ws.Range("E20").Value2 = -19.25730036044673
$ws.Range("F20").Value2 = 3.145986597812511
$ws.Range("G20").Value2 = -19.25730036044673
$ws.Range("H20").Value2 = 1.506713500706748
$ws.Range("I20").Value2 = 1.264990936851804
$ws.Range("J20").Value2 = -19.25730036044673
$ws.Range("K20").Value2 = 4.321925907763245
$ws.Range("B21").Value2 = -19.25730036044673
$ws.Range("C21").Value2 = 1.311480575288055
$ws.Range("D21").Value2 = -19.25730036044673
$ws.Range("E21").Value2 = 1.685979234021066
$ws.Range("F21").Value2 = -19.25730036044673
$ws.Range("G21").Value2 = 2.423888713880968
$ws.Range("H21").Value2 = 1.516435283806679
$ws.Range("I21").Value2 = -19.25730036044673
$ws.Range("J21").Value2 = -19.25730036044673
$ws.Range("K21").Value2 = -19.25730036044673
